$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.453.63"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "2.584.14"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "551.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.82%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.593"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.92%  "
$ws.Range("D9").Value = "2.596.82"
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.72"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.104"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.161"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.354"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.91%  "
$ws.Range("D14").Value = "3.043.03"
$ws.Range("E14").Value = "  -0.77%  "
$ws.Range("D15").Value = "59.449.32"
$ws.Range("E15").Value = "  -0.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.20"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000137"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.68%  "
$ws.Range("D18").Value = "2.586.70"
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "339.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.48"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.26%  "
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.477"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.32%  "
$ws.Range("E26").Value = "  -0.24%  "
$ws.Range("E27").Value = "  -1.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.45"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.18%  "
$ws.Range("D29").Value = "0.0₃0770"
$ws.Range("E29").Value = "  -1.31%  "
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.95%  "
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "158.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.16"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.13"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.74%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.17"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.902"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "37.58"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.02%  "
$ws.Range("B39").Value = "SuiNetwork"
$ws.Range("C39").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.842"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.16%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.47"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.66"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "287.80"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "135.54"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.82%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0972"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.597"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.46%  "
$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.68"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0533"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0234"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.06%  "
$ws.Range("D51").Value = "1.974.33"
$ws.Range("E51").Value = "  +1.78%  "
